# Fruta / hortaliza, semanal
# Insert two new price records (rows) right after the existing row 875
# ("Vega Modelo de Temuco" - Uva sheet), pushing all following rows down
# by two positions (old row 876 -> new row 878, ... old row 931 -> new row 933).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 876 (two single-row inserts keep the
# surrounding formatting - e.g. the date style in column D - consistent
# with the rest of the column, exactly like Excel does interactively).
$ws.Rows.Item(876).Insert()
$ws.Rows.Item(876).Insert()

# New row 876: Crimpson Seedless / Primera
$ws.Cells.Item(876, 1).Value2  = 10
$ws.Cells.Item(876, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(876, 3).Value2  = "La Araucanía"
$ws.Cells.Item(876, 4).Value2  = 44753
$ws.Cells.Item(876, 5).Value2  = 9
$ws.Cells.Item(876, 6).Value2  = "Fruta"
$ws.Cells.Item(876, 7).Value2  = 100109
$ws.Cells.Item(876, 8).Value2  = "Uva"
$ws.Cells.Item(876, 9).Value2  = 100109001
$ws.Cells.Item(876, 10).Value2 = "Uva"
$ws.Cells.Item(876, 11).Value2 = "Crimpson Seedless"
$ws.Cells.Item(876, 12).Value2 = "Primera"
$ws.Cells.Item(876, 13).Value2 = 380
$ws.Cells.Item(876, 14).Value2 = 7000
$ws.Cells.Item(876, 15).Value2 = 7000
$ws.Cells.Item(876, 16).Value2 = 7000
$ws.Cells.Item(876, 17).Value2 = "$/bandeja 8 kilos"
$ws.Cells.Item(876, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(876, 19).Value2 = 875
$ws.Cells.Item(876, 20).Value2 = 8

# New row 877: Red Globe / Primera
$ws.Cells.Item(877, 1).Value2  = 10
$ws.Cells.Item(877, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(877, 3).Value2  = "La Araucanía"
$ws.Cells.Item(877, 4).Value2  = 44753
$ws.Cells.Item(877, 5).Value2  = 9
$ws.Cells.Item(877, 6).Value2  = "Fruta"
$ws.Cells.Item(877, 7).Value2  = 100109
$ws.Cells.Item(877, 8).Value2  = "Uva"
$ws.Cells.Item(877, 9).Value2  = 100109001
$ws.Cells.Item(877, 10).Value2 = "Uva"
$ws.Cells.Item(877, 11).Value2 = "Red Globe"
$ws.Cells.Item(877, 12).Value2 = "Primera"
$ws.Cells.Item(877, 13).Value2 = 310
$ws.Cells.Item(877, 14).Value2 = 12000
$ws.Cells.Item(877, 15).Value2 = 12000
$ws.Cells.Item(877, 16).Value2 = 12000
$ws.Cells.Item(877, 17).Value2 = "$/bandeja 8 kilos"
$ws.Cells.Item(877, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(877, 19).Value2 = 1500
$ws.Cells.Item(877, 20).Value2 = 8
